$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Before / During / After" explainer rows (93-100) with their
# English source text in column A and Russian translation in column B,
# in the same row-major order they were authored so the shared-string
# table comes out in the same sequence.

$ws.Range("A93").Value = "Before"
$ws.Range("B93").Value = "Перед "

$ws.Range("A94").Value = "During"
$ws.Range("B94").Value = "В течение"

$ws.Range("A95").Value = "After"
$ws.Range("B95").Value = "После"

$ws.Range("A96").Value = "What's Happened?"
$ws.Range("B96").Value = "Что случилось?"

$ws.Range("A97").Value = "What's the Worst?"
$ws.Range("B97").Value = "Что хуже всего?"

$ws.Range("A98").Value = "Cascadia Quake"
$ws.Range("B98").Value = "Разлом Каскадия"

$ws.Range("A99").Value = "Tsunami Zone"
$ws.Range("B99").Value = "Зона цунами"

$ws.Range("A100").Value = "If the dams failed"
$ws.Range("B100").Value = "Если прорвало плотину"

# Match the author's final selection after typing the new block.
$ws.Range("A93:B100").Select() | Out-Null
